# "Penalty Reward System" edit: the forecast week-window rolled forward by
# one week (every Week_Start_Date becomes the following week's original
# date), the first three forecasted weeks got zeroed out (the "penalty"),
# and the Summary sheet's derived stats were updated to match.

$wb  = $excel.ActiveWorkbook

# --- Sheet "Forecast Comparison" ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# Week_Start_Date values are stored as plain text (e.g. "2025-01-12"), not
# real dates. Force the cells to Text format first so Excel doesn't
# reinterpret the date-shaped strings as date serials.
$ws1.Range("B2:B17").NumberFormat = "@"

$dates = @(
    "2025-01-12", "2025-01-19", "2025-01-26", "2025-02-02", "2025-02-09",
    "2025-02-16", "2025-02-23", "2025-03-02", "2025-03-09", "2025-03-16",
    "2025-03-23", "2025-03-30", "2025-04-06", "2025-04-13", "2025-04-20",
    "2025-04-27"
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 2).Value = $dates[$i]
}

# MyForecast for the (now) first three weeks is penalized down to 0.
$ws1.Cells.Item(2, 4).Value = 0
$ws1.Cells.Item(3, 4).Value = 0
$ws1.Cells.Item(4, 4).Value = 0

# --- Sheet "Summary" ---
$ws2 = $wb.Worksheets.Item("Summary")

# These cells also hold plain text (even the numeric-looking ones), so keep
# them as Text to avoid Excel auto-converting to numbers/dates.
$ws2.Range("B9:B13").NumberFormat = "@"
$ws2.Range("B15").NumberFormat = "@"

$ws2.Cells.Item(2, 2).Value  = "2023-01-01 to 2025-01-05"
$ws2.Cells.Item(9, 2).Value  = "1"
$ws2.Cells.Item(10, 2).Value = "1"
$ws2.Cells.Item(11, 2).Value = "1"
$ws2.Cells.Item(12, 2).Value = "0"
$ws2.Cells.Item(13, 2).Value = "2025-01-12"
$ws2.Cells.Item(15, 2).Value = "2025-02-02"
